$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new blank column before column N (14th column), shifting
# N->O, O->P, P->Q
$ws.Columns("N").Insert()

$ws.Range("S10").Select()
